$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "Петров"

$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("D12").Value = 1

$ws.Range("B21").Value = "5-"
$ws.Range("B16").Value = "4+"
$ws.Range("B18").Value = "4+"
$ws.Range("B17").Value = 5
$ws.Range("B19").Value = 4
$ws.Range("B20").Value = 4
$ws.Range("B22").Value = 5
$ws.Range("B23").Value = 5
$ws.Range("B24").Value = 5

$ws.Range("C11").Select()
